$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.437.17'
$ws.Range("E2").Value = '  +3.01%  '
$ws.Range("D3").Value = '2.992.13'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.76%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +2.19%  '
$ws.Range("D9").Value = '2.979.02'
$ws.Range("E9").Value = '  +2.24%  '
$ws.Range("E10").Value = '  +5.05%  '
$ws.Range("E11").Value = '  +11.23%  '
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.63%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '3.489.43'
$ws.Range("E16").Value = '  +2.45%  '
$ws.Range("E17").Value = '  +4.75%  '
$ws.Range("D18").Value = '2.990.64'
$ws.Range("E18").Value = '  +2.42%  '
$ws.Range("D19").Value = '59.427.20'
$ws.Range("E19").Value = '  +2.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '434.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.16%  '
$ws.Range("E22").Value = '  +4.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").Value = '  +11.16%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +3.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.76'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.76%  '
$ws.Range("E31").Value = '  +5.25%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.106'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.68%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.72'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("D34").Value = '0.0₃0776'
$ws.Range("E34").Value = '  +12.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.990'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.71%  '
$ws.Range("E36").Value = '  +4.14%  '
$ws.Range("E37").Value = '  +1.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("E39").Value = '  -1.40%  '
$ws.Range("E40").Value = '  +7.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '400.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0353'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.87%  '
$ws.Range("D43").Value = '2.759.97'
$ws.Range("E43").Value = '  +2.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.105'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  +7.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +18.97%  '
$ws.Range("E50").Value = '  +3.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.58%  '
